$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").Value = "https://pubmed.ncbi.nlm.nih.gov/39226116/"
Write-Output $ws.Hyperlinks.Count()
